$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "29.279.74", "  -0.36%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.863.59", "  -0.54%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.001", "  +0.04%  "),
    @(5, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.7159", "  +0.43%  "),
    @(6, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "240.85", "  +0.55%  "),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  +0.01%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3095", "  +0.50%  "),
    @(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07731", "  -0.58%  "),
    @(10, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "25.16", "  +0.75%  "),
    @(11, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.263.83", "  +20.69%  "),
    @(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.08319", "  +0.96%  "),
    @(13, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.7190", "  -0.66%  "),
    @(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.230", "  -0.69%  "),
    @(15, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "90.81", "  -0.72%  "),
    @(16, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "29.282.67", "  -0.28%  "),
    @(17, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.921", "  +0.70%  "),
    @(18, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "243.78", "  +0.73%  "),
    @(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007821", "  -1.18%  "),
    @(20, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.120.36", "  +0.15%  "),
    @(21, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "13.17", "  -0.90%  "),
    @(22, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  +0.13%  "),
    @(23, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "8.020", "  +1.66%  "),
    @(24, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  +0.05%  "),
    @(25, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1633", "  +5.10%  "),
    @(26, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "162.63", "  -0.56%  "),
    @(27, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.936", "  -0.77%  "),
    @(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.72", "  +2.33%  "),
    @(29, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.351", "  -0.43%  "),
    @(30, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.496", "  +1.21%  "),
    @(31, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.439", "  +1.72%  "),
    @(32, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.247", "  +3.39%  "),
    @(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05204", "  -1.18%  "),
    @(34, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.930", "  +0.26%  "),
    @(35, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.174", "  -1.92%  "),
    @(36, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7283", "  +1.72%  "),
    @(37, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.685", "  +0.29%  "),
    @(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01863", "  +0.26%  "),
    @(39, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.700", "  -0.41%  "),
    @(40, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.185.14", "  -0.59%  "),
    @(41, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9091", "  +0.34%  "),
    @(42, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.202", "  +2.13%  "),
    @(43, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "72.98", "  +1.12%  "),
    @(44, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9998", "  -0.08%  "),
    @(45, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "103.05", "  +0.33%  "),
    @(46, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "2.018.27", "  +0.71%  "),
    @(47, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.5225", "  -2.31%  "),
    @(48, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.777", "  +1.31%  "),
    @(49, "SynthetixNetwork", "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx", "2.908", "  +0.21%  "),
    @(50, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.356", "  +1.22%  "),
    @(51, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.096", "  +1.27%  "),
)

foreach ($item in $data) {
    $r = $item[0]
    foreach ($col in 2..5) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $item[$col - 1]
        $cell.Style = "Normal"
    }
}

$wb.Save()
